# Updated cryptos list (price / 1h-volume-change refresh, plus the
# EnergySwap/Filecoin row swap) to match the latest GitHub Actions scrape.
#
# Price cells (column D) that are purely numeric-looking strings would be
# auto-coerced to Number by a plain .Value assignment (e.g. "44.20" -> 44.2,
# losing the trailing zero, or "0.0000307" -> scientific notation). To keep
# them as literal text (matching the original inlineStr cells) we briefly
# force NumberFormat "@" before writing, then restore the cell style to
# "Normal" so the cell's style reference is left exactly as it was.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.500.43"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.708.96"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.47%  "
$ws.Range("E6").Value = "  -2.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "651.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.34%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  -6.92%  "
$ws.Range("D11").Value = "3.709.11"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000307"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +14.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "44.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "4.399.53"
$ws.Range("E16").Value = "  +0.09%  "
$ws.Range("D17").Value = "96.263.00"
$ws.Range("E17").Value = "  -0.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "3.707.82"
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -8.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "520.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000211"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "101.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.178"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.88"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.996"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.28"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "646.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.67%  "
$ws.Range("E38").Value = "  -1.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.46%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "41.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.72%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.22%  "
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0453"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.78%  "
$ws.Range("E51").Value = "  +1.15%  "
